# [ISP-2] Configurazione forecast per Intesa
#
# Adds RPM / intesa_user visibility configuration to the
# "Page Group Visibility" sheet, updates the active sheet / selections
# accordingly, matching the author's edit.

$wb = $excel.ActiveWorkbook

# --- Permissions sheet: just a selection change (no data change) ---
$wsPermissions = $wb.Worksheets.Item("Permissions")
$wsPermissions.Range("A3").Select() | Out-Null

# --- Page Group Visibility sheet: add the new RPM rows ---
$wsPageGroupVisibility = $wb.Worksheets.Item("Page Group Visibility")

$wsPageGroupVisibility.Range("A2").Value = "CREATE/MODIFY"
$wsPageGroupVisibility.Range("B2").Value = "RPM"
$wsPageGroupVisibility.Range("C2").Value = "core_group_admin"

$wsPageGroupVisibility.Range("A3").Value = "CREATE/MODIFY"
$wsPageGroupVisibility.Range("B3").Value = "RPM"
$wsPageGroupVisibility.Range("C3").Value = "intesa_user"

# Give C3 its own (otherwise identical) style, as in the source workbook
$wsPageGroupVisibility.Range("C3").Font.Name = "Trebuchet MS"
$wsPageGroupVisibility.Range("C3").Font.Size = 10

# --- Activities Labels sheet: just a selection change ---
$wsActivitiesLabels = $wb.Worksheets.Item("Activities Labels")
$wsActivitiesLabels.Range("B7").Select() | Out-Null

# --- Make "Page Group Visibility" the active sheet/selection last, ---
# --- so it ends up tabSelected + becomes the active tab of the book ---
$wsPageGroupVisibility.Range("A3:B3").Select() | Out-Null
